# Generate Report for Handback
#
# The handback transform failed for the second tracked file
# (7cefde0a-1ef2-4531-bc1f-fca4fc8bf6f8.md) in both the zh-cn and de-de
# locales: the handback file name didn't match the expected handoff file
# name. Update the status for that row and record the error detail.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

$newStatus = "Handback transform failed"

# Status column (C) for the second file row (row 3) flips from
# "Ready for handoff" to "Handback transform failed" in both locale
# sheets, and the Overview rollup (columns B/C) mirrors it.
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# Record the handback/handoff file name mismatch in the Error Detail
# column (K) for row 3 on each locale sheet.
$zhcn.Range("K3").Value = "Handback file name: 505oeiuv.whq is different with handoff file name: 7cefde0a-1ef2-4531-bc1f-fca4fc8bf6f8.3b7f5b2e229bf6c14e80359c3cca8c3f45084720.zh-cn."
$dede.Range("K3").Value = "Handback file name: 505oeiuv.whq is different with handoff file name: 7cefde0a-1ef2-4531-bc1f-fca4fc8bf6f8.3b7f5b2e229bf6c14e80359c3cca8c3f45084720.de-de."
